# Applies the "Updated symbol list" crypto price/volume/coin refresh described
# in the commit diff. Column B/C (Coin, Link) cells get plain string writes;
# column D/E (Price, Volume) cells look numeric/percentage-like, so Excel would
# normally auto-convert them to Number/Percentage on assignment. The source
# workbook stores them as literal text (inlineStr), so we force text storage
# (NumberFormat "@") for the assignment and then restore the "Normal" style so
# no stray per-cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link columns (plain text, no numeric coercion risk) ---
$textCells = @(
    @{Cell="B7"; Value="GateToken"},
    @{Cell="C7"; Value="https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"},
    @{Cell="B8"; Value="MXToken"},
    @{Cell="C8"; Value="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"},
    @{Cell="B9"; Value="LiechtensteinCryptoassetsExchange"},
    @{Cell="C9"; Value="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"},
    @{Cell="B10"; Value="WazirX"},
    @{Cell="C10"; Value="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"},
    @{Cell="B11"; Value="MCDex"},
    @{Cell="C11"; Value="https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"},
    @{Cell="B12"; Value="MandalaExchangeToken"},
    @{Cell="C12"; Value="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"},
    @{Cell="B13"; Value="BitrueCoin"},
    @{Cell="C13"; Value="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"},
    @{Cell="B14"; Value="BitMartToken"},
    @{Cell="C14"; Value="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"},
    @{Cell="B15"; Value="BitForexToken"},
    @{Cell="C15"; Value="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"},
    @{Cell="B16"; Value="CoinExToken"},
    @{Cell="C16"; Value="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"},
    @{Cell="B17"; Value="TigerCash"},
    @{Cell="C17"; Value="https://coinranking.com/coin/6hIn06L2+tigercash-tch"},
    @{Cell="B18"; Value="LEO"},
    @{Cell="C18"; Value="https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"}
)

foreach ($item in $textCells) {
    $ws.Range($item.Cell).Value = $item.Value
}

# --- Price / Volume(1h) columns (force text so "329.12" / "0.21%" stay literal) ---
$numericLookingCells = @(
    @{Cell="D2"; Value="329.12"},
    @{Cell="E2"; Value="0.21%"},
    @{Cell="D3"; Value="44.10"},
    @{Cell="E3"; Value="-0.05%"},
    @{Cell="D4"; Value="5.511"},
    @{Cell="E4"; Value="-1.19%"},
    @{Cell="D5"; Value="0.08055"},
    @{Cell="E5"; Value="-0.14%"},
    @{Cell="D6"; Value="2.026"},
    @{Cell="E6"; Value="6.16%"},
    @{Cell="D7"; Value="4.410"},
    @{Cell="E7"; Value="2.99%"},
    @{Cell="D8"; Value="0.9558"},
    @{Cell="E8"; Value="0.33%"},
    @{Cell="D9"; Value="0.1130"},
    @{Cell="E9"; Value="-7.15%"},
    @{Cell="D10"; Value="0.1872"},
    @{Cell="E10"; Value="1.42%"},
    @{Cell="D11"; Value="10.09"},
    @{Cell="E11"; Value="1.03%"},
    @{Cell="D12"; Value="0.1001"},
    @{Cell="E12"; Value="3.48%"},
    @{Cell="D13"; Value="0.04837"},
    @{Cell="E13"; Value="9.70%"},
    @{Cell="D14"; Value="0.1058"},
    @{Cell="E14"; Value="-0.77%"},
    @{Cell="D15"; Value="0.001257"},
    @{Cell="E15"; Value="-2.97%"},
    @{Cell="D16"; Value="0.04090"},
    @{Cell="E16"; Value="-2.60%"},
    @{Cell="D17"; Value="0.006052"},
    @{Cell="E17"; Value="2.00%"},
    @{Cell="D18"; Value="3.370"},
    @{Cell="E18"; Value="-0.68%"},
    @{Cell="D19"; Value="2.602"},
    @{Cell="E19"; Value="1.50%"},
    @{Cell="D20"; Value="0.3283"},
    @{Cell="E20"; Value="-4.29%"},
    @{Cell="D21"; Value="0.1400"},
    @{Cell="E21"; Value="-0.82%"},
    @{Cell="D22"; Value="0.2571"},
    @{Cell="E22"; Value="2.77%"},
    @{Cell="D23"; Value="0.001306"},
    @{Cell="E23"; Value="5.09%"},
    @{Cell="D24"; Value="0.004360"},
    @{Cell="E24"; Value="0.55%"},
    @{Cell="D25"; Value="0.0001250"},
    @{Cell="E25"; Value="5.14%"},
    @{Cell="D26"; Value="0.0003738"},
    @{Cell="E26"; Value="-6.04%"},
    @{Cell="D38"; Value="0.02591"},
    @{Cell="E38"; Value="-3.10%"},
    @{Cell="D39"; Value="0.05671"},
    @{Cell="E39"; Value="2.55%"},
    @{Cell="D40"; Value="0.007656"},
    @{Cell="E40"; Value="0.53%"},
    @{Cell="D41"; Value="0.1403"},
    @{Cell="E41"; Value="-0.37%"},
    @{Cell="D42"; Value="0.007369"},
    @{Cell="E42"; Value="-6.35%"},
    @{Cell="D43"; Value="0.002007"},
    @{Cell="E43"; Value="-0.29%"},
    @{Cell="D44"; Value="0.008520"},
    @{Cell="E44"; Value="-4.19%"},
    @{Cell="D45"; Value="0.00007077"},
    @{Cell="E45"; Value="1.80%"},
    @{Cell="D46"; Value="0.00000000750"},
    @{Cell="E46"; Value="0.08%"},
    @{Cell="D47"; Value="0.0005802"},
    @{Cell="E47"; Value="-0.17%"},
    @{Cell="D48"; Value="0.003497"},
    @{Cell="E48"; Value="54.22%"},
    @{Cell="D49"; Value="0.003506"},
    @{Cell="E49"; Value="26.31%"},
    @{Cell="D50"; Value="0.00002099"},
    @{Cell="E50"; Value="0.08%"},
    @{Cell="D51"; Value="0.0001999"},
    @{Cell="E51"; Value="0.08%"}
)

foreach ($item in $numericLookingCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}
